$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2017-12-31 00:00:00"
$ws.Range("O2").Value = 14058996.28
$ws.Range("P2").Value = 14.8292509966
$ws.Range("Q2").Value = 1508449549.66
$ws.Range("R2").Value = 1591.0934566102
$ws.Range("S2").Value = 419925817.06
$ws.Range("T2").Value = 442.9324268329
$ws.Range("U2").Value = -44203250.75
$ws.Range("V2").Value = -46.6250283578
$ws.Range("Y2").Value = 44278869.86
$ws.Range("Z2").Value = 46.7047904361
$ws.Range("AA2").Value = 127199771.53
$ws.Range("AB2").Value = 134.168706013
$ws.Range("AC2").Value = 94805842.06999999
$ws.Range("AD2").Value = 114.9012626191
